$wb = $excel.ActiveWorkbook

# --- Sheet 1: Metadata ---
$meta = $wb.Worksheets.Item("Metadata")

# Version: 1.1.0 -> 2.0.0
$meta.Range("B3").Value = "2.0.0"

# Date: 2023-07-10T23:08:03+02:00 -> 2024-06-04T14:59:10+02:00
$meta.Range("B8").Value = "2024-06-04T14:59:10+02:00"

# Contact: "No display for ContactDetail" -> "Kommunernes Landsforening (http://kl.dk)"
$meta.Range("B10").Value = "Kommunernes Landsforening (http://kl.dk)"

# --- Sheet 2: Include from FSIII ---
$ws = $wb.Worksheets.Item("Include from FSIII")

# Insert 12 new blank rows above the existing concept rows (rows 2-13),
# pushing the current content down to rows 14-27.
$ws.Range("A2:A13").Insert()

# Copy the formatting (style) of the row just below (now row 14, style "s=2")
# onto the newly inserted blank rows so they match the rest of the table.
$ws.Range("A14:B14").Copy()
$ws.Range("A2:B13").PasteSpecial(-4122)

# Fill in column A of the newly inserted rows with the new concept identifiers.
$newConcepts = @(
  "94e9c867-fbc8-4d35-8596-e6b8765b12e8",
  "55670b1e-7a36-46b2-8712-b7536237f22d",
  "9162d29a-1c7f-4585-8145-8fb4f1a999e3",
  "fa6aa904-d06e-4029-b4c4-13ead04ace27",
  "3f00a76f-8e7b-4b13-80cc-f2ceef4e51d1",
  "01150cdb-6098-48ce-bb61-60967f6bcc37",
  "1bb534f3-e526-41a9-b9c3-6157ea19c915",
  "cc377732-7f14-49b7-8940-1aa07b8884e7",
  "25dcedb3-7149-4ef9-a2c3-be30267441fb",
  "045fa500-35b0-46b7-97dd-adb60888a8ea",
  "8c539fd9-7f31-4b4e-8b30-8298c8ab640f",
  "5bfe4bda-2358-41da-946e-1fdaa33d5fe8"
)

for ($i = 0; $i -lt $newConcepts.Length; $i++) {
  $row = 2 + $i
  $ws.Cells.Item($row, 1).Value = $newConcepts[$i]
}
